$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of per-emotion offer/self-report data for 5 additional videos
# (1032, 1033, 1035, 1039, 1038), inserted above the trailing "NULL" sentinel row.
$newData = @(
    @(1032, "angry", "1.0359869138494999E-2"),
    @(1032, "disgust", "0"),
    @(1032, "fear", "0"),
    @(1032, "happy", "0"),
    @(1032, "sad", "1.30861504907306E-2"),
    @(1032, "surprise", "0"),
    @(1032, "neutral", "0.97655398037077401"),
    @(1033, "angry", "2.4840312278211499E-2"),
    @(1033, "disgust", "0"),
    @(1033, "fear", "7.0972320794889996E-4"),
    @(1033, "happy", "0.49609652235628099"),
    @(1033, "sad", "0.258339247693399"),
    @(1033, "surprise", "0"),
    @(1033, "neutral", "0.22001419446415801"),
    @(1035, "angry", "0.224587315377932"),
    @(1035, "disgust", "0"),
    @(1035, "fear", "0"),
    @(1035, "happy", "0"),
    @(1035, "sad", "0.719374456993918"),
    @(1035, "surprise", "0"),
    @(1035, "neutral", "5.6038227628149397E-2"),
    @(1039, "angry", "0"),
    @(1039, "disgust", "0"),
    @(1039, "fear", "0"),
    @(1039, "happy", "6.8645640074211506E-2"),
    @(1039, "sad", "2.7829313543599201E-3"),
    @(1039, "surprise", "0"),
    @(1039, "neutral", "0.92857142857142805"),
    @(1038, "angry", "2.9970029970029901E-3"),
    @(1038, "disgust", "0"),
    @(1038, "fear", "9.99000999000999E-4"),
    @(1038, "happy", "2.9970029970029899E-2"),
    @(1038, "sad", "0.84415584415584399"),
    @(1038, "surprise", "0"),
    @(1038, "neutral", "0.121878121878121")
)

# Insert 35 blank rows at row 218, pushing the existing row 218 (the NULL
# sentinel row) down to row 253, copying formatting from the row above.
$insertRange = $ws.Range("A218:A" + (218 + $newData.Count - 1))
$insertRange.EntireRow.Insert() | Out-Null

$startRow = 218
$i = 0
foreach ($item in $newData) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $r
    $ws.Cells.Item($r, 2).Value = $item[0]
    $ws.Cells.Item($r, 3).Value = $item[1]
    $ws.Cells.Item($r, 4).Value = [double]$item[2]
    $i = $i + 1
}

# Match the author's final cursor position.
$ws.Range("E9").Select() | Out-Null
